# B1--and-B2-PowerPoint.pptx -- "Tue, Jul 07, 2020 12:05:36 PM" commit
#
# Two logical changes are made:
#   1. The table on slide 5 gets a different built-in table style applied
#      (its a:tableStyleId GUID changes).
#   2. The deck's theme is swapped from the "Integral" (Red Violet) colour
#      palette to the stock "Office Theme" colour palette, i.e. each of the
#      12 theme colour slots is updated to the standard Office RGB values.

function ConvertTo-PpRGB([string]$hex) {
    # PowerPoint/VBA RGB values are packed as 0x00BBGGRR (low byte = Red).
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 --------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{1AED6874-C10F-4AAF-A7D5-7690D1565699}")

# --- 2. Re-colour the theme to the stock "Office Theme" palette ------------
# ThemeColorScheme.Item index order is: dk1, lt1, dk2, lt2,
# accent1-accent6, hlink, folHlink.
$officeThemeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$colorScheme = $p.Designs.Item(1).SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-PpRGB $officeThemeHex[$i - 1]
}
